# Added total heatmap generation:
# Insert a new column A in front of the existing data. The former columns
# A:E (startHeight, stopHeight, startWidth, stopWidth, Name) shift right to
# become B:F. The new column A holds a zero-based row index (0, 1, 2, ...)
# for each data row, formatted with the same (bold/centered/bordered)
# header style used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:E one column to the right, creating new column A.
$ws.Columns.Item(1).EntireColumn.Insert()

# Determine how many data rows exist below the header row (row 1) using the
# data now sitting in column B (the original column A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Fill new column A (rows 2..lastRow) with a 0-based running index.
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

# Apply the header cell's style (bold, centered, bordered) to the new index
# column's data cells, matching the style used by row 1's header cells.
$ws.Range("B1").Copy()
$ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1)).PasteSpecial(-4122)
